$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain text (matches original inlineStr text cells)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '41.151.54'
$ws.Range("E2").Value = '  -3.56%  '
$ws.Range("D3").Value = '2.461.11'
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '311.92'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").Value = '93.80'
$ws.Range("E6").Value = '  -6.14%  '
$ws.Range("E7").Value = '  -2.94%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -4.71%  '
$ws.Range("D10").Value = '33.14'
$ws.Range("E10").Value = '  -6.96%  '
$ws.Range("E11").Value = '  -3.29%  '
$ws.Range("E12").Value = '  -1.37%  '
$ws.Range("D13").Value = '6.96'
$ws.Range("E13").Value = '  -4.80%  '
$ws.Range("D14").Value = '2.840.86'
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("D15").Value = '2.463.34'
$ws.Range("E15").Value = '  -4.17%  '
$ws.Range("D16").Value = '14.80'
$ws.Range("E16").Value = '  -2.94%  '
$ws.Range("E17").Value = '  -3.45%  '
$ws.Range("D18").Value = '41.122.65'
$ws.Range("E18").Value = '  -3.59%  '
$ws.Range("D19").Value = '6.27'
$ws.Range("E19").Value = '  -5.98%  '
$ws.Range("D20").Value = '0.0₃0919'
$ws.Range("E20").Value = '  -3.07%  '
$ws.Range("D21").Value = '11.12'
$ws.Range("E21").Value = '  -8.80%  '
$ws.Range("D22").Value = '68.44'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("D23").Value = '235.25'
$ws.Range("E23").Value = '  -2.98%  '
$ws.Range("E24").Value = '  -3.90%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -6.04%  '
$ws.Range("D27").Value = '23.95'
$ws.Range("E27").Value = '  -5.79%  '
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  -6.23%  '
$ws.Range("D29").Value = '9.58'
$ws.Range("E29").Value = '  -5.52%  '
$ws.Range("D30").Value = '36.20'
$ws.Range("E30").Value = '  -5.70%  '
$ws.Range("D31").Value = '152.54'
$ws.Range("E31").Value = '  -4.03%  '
$ws.Range("E32").Value = '  -5.15%  '
$ws.Range("E33").Value = '  -5.43%  '
$ws.Range("D34").Value = '2.55'
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("D35").Value = '0.0741'
$ws.Range("E35").Value = '  -5.03%  '
$ws.Range("D36").Value = '3.01'
$ws.Range("E36").Value = '  -3.36%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").Value = '17.01'
$ws.Range("E37").Value = '  -7.58%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '1.87'
$ws.Range("E38").Value = '  -4.08%  '
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("E40").Value = '  -8.07%  '
$ws.Range("D41").Value = '4.19'
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E43").Value = '  -11.04%  '
$ws.Range("D44").Value = '1.971.44'
$ws.Range("E45").Value = '  -5.13%  '
$ws.Range("D46").Value = '3.03'
$ws.Range("E46").Value = '  -7.79%  '
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '96.77'
$ws.Range("E48").Value = '  -3.89%  '
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = '68.85'
$ws.Range("E49").Value = '  -4.27%  '
$ws.Range("D50").Value = '0.177'
$ws.Range("E50").Value = '  -6.32%  '
$ws.Range("D51").Value = '73.68'
$ws.Range("E51").Value = '  -6.95%  '
